$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 58:59, shifting the existing rows 58.. down to 60..
$ws.Rows("58:59").Insert()

# Populate the two new rows with the new "Patagonia" / "1a (cosecha)" record
# (the same record is duplicated on both row 58 and row 59).
foreach ($r in 58, 59) {
    $ws.Cells.Item($r, 1).Value = 1
    $ws.Cells.Item($r, 2).Value = "Agrícola del Norte S.A. de Arica"
    $ws.Cells.Item($r, 3).Value = "Arica y Parinacota"
    $ws.Cells.Item($r, 4).Value = 44923
    $ws.Cells.Item($r, 5).Value = 15
    $ws.Cells.Item($r, 6).Value = 100114001
    $ws.Cells.Item($r, 7).Value = "Papa"
    $ws.Cells.Item($r, 8).Value = "Patagonia"
    $ws.Cells.Item($r, 9).Value = "1a (cosecha)"
    $ws.Cells.Item($r, 10).Value = 1000
    $ws.Cells.Item($r, 11).Value = 14500
    $ws.Cells.Item($r, 12).Value = 15000
    $ws.Cells.Item($r, 13).Value = 14800
    $ws.Cells.Item($r, 14).Value = "$/saco 25 kilos"
    $ws.Cells.Item($r, 15).Value = "Región del Maule"
    $ws.Cells.Item($r, 16).Value = 592
    $ws.Cells.Item($r, 17).Value = 25
    $ws.Cells.Item($r, 18).Value = "Hortaliza"
}
